$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 2
Set-TextValue $ws.Range("D2") '43.123.35'
Set-TextValue $ws.Range("E2") '  +2.20%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.556.34'
Set-TextValue $ws.Range("E3") '  +1.89%  '

# Row 4
Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  -0.13%  '

# Row 5
Set-TextValue $ws.Range("D5") '316.91'
Set-TextValue $ws.Range("E5") '  +0.91%  '

# Row 6
Set-TextValue $ws.Range("D6") '97.46'
Set-TextValue $ws.Range("E6") '  +4.45%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +1.03%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.03%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.541'
Set-TextValue $ws.Range("E9") '  +3.34%  '

# Row 10
Set-TextValue $ws.Range("D10") '35.71'
Set-TextValue $ws.Range("E10") '  +1.30%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0812'
Set-TextValue $ws.Range("E11") '  +1.43%  '

# Row 12
Set-TextValue $ws.Range("D12") '7.51'
Set-TextValue $ws.Range("E12") '  +1.38%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D13") '0.108'
Set-TextValue $ws.Range("E13") '  -4.52%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D14") '2.948.13'
Set-TextValue $ws.Range("E14") '  +1.60%  '

# Row 15
Set-TextValue $ws.Range("D15") '2.547.14'
Set-TextValue $ws.Range("E15") '  +1.45%  '

# Row 16
Set-TextValue $ws.Range("E16") '  -1.22%  '

# Row 17
Set-TextValue $ws.Range("D17") '0.851'
Set-TextValue $ws.Range("E17") '  +1.68%  '

# Row 18
Set-TextValue $ws.Range("D18") '43.125.31'
Set-TextValue $ws.Range("E18") '  +1.84%  '

# Row 19
Set-TextValue $ws.Range("E19") '  +5.36%  '

# Row 20
Set-TextValue $ws.Range("D20") '12.61'
Set-TextValue $ws.Range("E20") '  -1.10%  '

# Row 21
Set-TextValue $ws.Range("D21") '0.0₃0965'
Set-TextValue $ws.Range("E21") '  +1.52%  '

# Row 22
Set-TextValue $ws.Range("D22") '70.09'
Set-TextValue $ws.Range("E22") '  -0.29%  '

# Row 23
Set-TextValue $ws.Range("D23") '253.58'
Set-TextValue $ws.Range("E23") '  +2.06%  '

# Row 24
Set-TextValue $ws.Range("E24") '  +1.50%  '

# Row 25
Set-TextValue $ws.Range("E25") '  +3.35%  '

# Row 26
Set-TextValue $ws.Range("E26") '  +3.16%  '

# Row 27
Set-TextValue $ws.Range("E27") '  +0.26%  '

# Row 28
Set-TextValue $ws.Range("D28") '2.44'
Set-TextValue $ws.Range("E28") '  +3.41%  '

# Row 29
Set-TextValue $ws.Range("D29") '40.92'
Set-TextValue $ws.Range("E29") '  +6.06%  '

# Row 30
Set-TextValue $ws.Range("D30") '10.33'
Set-TextValue $ws.Range("E30") '  +2.56%  '

# Row 31
Set-TextValue $ws.Range("E31") '  +0.42%  '

# Row 32
Set-TextValue $ws.Range("E32") '  -0.30%  '

# Row 33
Set-TextValue $ws.Range("D33") '19.32'
Set-TextValue $ws.Range("E33") '  -0.74%  '

# Row 34
Set-TextValue $ws.Range("E34") '  +3.51%  '

# Row 35
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D35") '2.12'
Set-TextValue $ws.Range("E35") '  +1.17%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D36") '0.0800'
Set-TextValue $ws.Range("E36") '  +3.17%  '

# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D37") '3.32'
Set-TextValue $ws.Range("E37") '  +1.89%  '

# Row 38
Set-TextValue $ws.Range("E38") '  +2.65%  '

# Row 39
Set-TextValue $ws.Range("E39") '  +4.85%  '

# Row 40
Set-TextValue $ws.Range("E40") '  +0.61%  '

# Row 41
Set-TextValue $ws.Range("D41") '22.02'
Set-TextValue $ws.Range("E41") '  -6.48%  '

# Row 42
Set-TextValue $ws.Range("D42") '3.85'
Set-TextValue $ws.Range("E42") '  +2.76%  '

# Row 43
Set-TextValue $ws.Range("E43") '  +2.41%  '

# Row 44
Set-TextValue $ws.Range("E44") '  -0.08%  '

# Row 45
Set-TextValue $ws.Range("E45") '  -0.32%  '

# Row 46
Set-TextValue $ws.Range("D46") '1.989.80'
Set-TextValue $ws.Range("E46") '  -0.85%  '

# Row 47
Set-TextValue $ws.Range("D47") '9.11'
Set-TextValue $ws.Range("E47") '  +3.99%  '

# Row 48
Set-TextValue $ws.Range("D48") '84.90'
Set-TextValue $ws.Range("E48") '  +1.32%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.804.59'
Set-TextValue $ws.Range("E49") '  +1.64%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D50") '104.85'
Set-TextValue $ws.Range("E50") '  +3.52%  '

# Row 51
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue $ws.Range("D51") '74.44'
Set-TextValue $ws.Range("E51") '  +3.06%  '
